$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds the "Förändrad" date, currently 45171 (2023-09-02) for all
# data rows; update it to 45172 (2023-09-03) for rows 2 through 44.
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 3).Value = 45172
}
